$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 66.23077000000001
$ws.Range("I11").Value = 66.23077000000001
$ws.Range("K11").Value = 66.23077000000001
$ws.Range("M11").Value = 73.76922999999999
$ws.Range("H32").Value = 5134.8667
$ws.Range("I32").Value = 2595.6667
$ws.Range("J32").Value = 5769.6665
$ws.Range("K32").Value = 2595.6667
$ws.Range("L32").Value = 5769.6665
$ws.Range("M32").Value = -2269.6667
$ws.Range("N32").Value = -6421.6665
$ws.Range("H51").Value = 9479
$ws.Range("J51").Value = 9759.781999999999
$ws.Range("L51").Value = 9759.781999999999
$ws.Range("N51").Value = -10727.782
$ws.Range("H92").Value = 50000220
$ws.Range("I92").Value = 66666800
$ws.Range("K92").Value = 66666800
$ws.Range("M92").Value = -66665552
$ws.Range("H96").Value = 797.25
$ws.Range("I96").Value = 899.7143
$ws.Range("K96").Value = 2699.1429
$ws.Range("M96").Value = -1326.1429
$ws.Range("H112").Value = 4108.154
$ws.Range("J112").Value = 4304.6665
$ws.Range("L112").Value = 12913.9995
$ws.Range("N112").Value = -15129.9995
$ws.Range("H116").Value = 19349524
$ws.Range("I116").Value = 38695440
$ws.Range("J116").Value = 3610.2856
$ws.Range("K116").Value = 38695440
$ws.Range("L116").Value = 3610.2856
$ws.Range("M116").Value = -38691998
$ws.Range("N116").Value = -10494.2856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 44910
$ws.Range("I32").Value = 48143.125
$ws.Range("K32").Value = 48143.125
$ws.Range("M32").Value = -47856.125
$ws.Range("H45").Value = 2922.4546
$ws.Range("I45").Value = 1945.5
$ws.Range("K45").Value = 1945.5
$ws.Range("M45").Value = -1568.5
$ws.Range("H61").Value = 3828.6287
$ws.Range("I61").Value = 2943.88
$ws.Range("K61").Value = 2943.88
$ws.Range("M61").Value = -2731.88
$ws.Range("H110").Value = 1201790.5
$ws.Range("I110").Value = 1361662.9
$ws.Range("J110").Value = 2747.5
$ws.Range("K110").Value = 1361662.9
$ws.Range("L110").Value = 2747.5
$ws.Range("M110").Value = -1359617.9
$ws.Range("N110").Value = -6837.5
$ws.Range("H122").Value = 4507417.5
$ws.Range("I122").Value = 4904675
$ws.Range("J122").Value = 5166.6665
$ws.Range("K122").Value = 14714025
$ws.Range("L122").Value = 15499.9995
$ws.Range("M122").Value = -14711575
$ws.Range("N122").Value = -20399.9995
$ws.Range("H132").Value = 16828.533
$ws.Range("I132").Value = 20588.621
$ws.Range("K132").Value = 61765.863
$ws.Range("M132").Value = -59235.863
$ws.Range("H136").Value = 3828.6287
$ws.Range("I136").Value = 2943.88
$ws.Range("K136").Value = 8831.639999999999
$ws.Range("M136").Value = -6281.639999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7492.8335
$ws.Range("I20").Value = 9710
$ws.Range("J20").Value = 5275.6665
$ws.Range("K20").Value = 9710
$ws.Range("L20").Value = 5275.6665
$ws.Range("M20").Value = -9463
$ws.Range("N20").Value = -5769.6665
$ws.Range("H94").Value = 14269936
$ws.Range("I94").Value = 17714106
$ws.Range("J94").Value = 1229.1428
$ws.Range("K94").Value = 17714106
$ws.Range("L94").Value = 1229.1428
$ws.Range("M94").Value = -17713655
$ws.Range("N94").Value = -2131.1428
$ws.Range("H134").Value = 2150.853
$ws.Range("I134").Value = 1554.7241
$ws.Range("K134").Value = 4664.1723
$ws.Range("M134").Value = -2129.1723

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2678.4243
$ws.Range("I7").Value = 2481.5652
$ws.Range("J7").Value = 3131.2
$ws.Range("K7").Value = 2481.5652
$ws.Range("L7").Value = 3131.2
$ws.Range("M7").Value = -2368.5652
$ws.Range("N7").Value = -3357.2
$ws.Range("H31").Value = 13515271
$ws.Range("I31").Value = 14707215
$ws.Range("J31").Value = 6573
$ws.Range("K31").Value = 14707215
$ws.Range("L31").Value = 6573
$ws.Range("M31").Value = -14706920
$ws.Range("N31").Value = -7163
$ws.Range("H34").Value = 13515271
$ws.Range("I34").Value = 14707215
$ws.Range("J34").Value = 6573
$ws.Range("K34").Value = 14707215
$ws.Range("L34").Value = 6573
$ws.Range("M34").Value = -14707013
$ws.Range("N34").Value = -6977
$ws.Range("H107").Value = 957974.9399999999
$ws.Range("J107").Value = 1554.3334
$ws.Range("L107").Value = 1554.3334
$ws.Range("N107").Value = -5394.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 34675320
$ws.Range("I4").Value = 40150990
$ws.Range("K4").Value = 120452970
$ws.Range("M4").Value = -120452858
$ws.Range("H33").Value = 435.77777
$ws.Range("I33").Value = 477.75
$ws.Range("K33").Value = 2866.5
$ws.Range("M33").Value = -2583.5
$ws.Range("H48").Value = 2198
$ws.Range("J48").Value = 2198
$ws.Range("L48").Value = 6594
$ws.Range("N48").Value = -7094
$ws.Range("H107").Value = 616
$ws.Range("J107").Value = 982
$ws.Range("L107").Value = 2946
$ws.Range("N107").Value = -6786
$ws.Range("H114").Value = 2319.2307
$ws.Range("I114").Value = 833.3333
$ws.Range("J114").Value = 2765
$ws.Range("K114").Value = 2499.9999
$ws.Range("L114").Value = 8295
$ws.Range("M114").Value = 754.0001000000002
$ws.Range("N114").Value = -14803
$ws.Range("H129").Value = 1260.1072
$ws.Range("I129").Value = 832.26086
$ws.Range("J129").Value = 3228.2
$ws.Range("K129").Value = 2496.78258
$ws.Range("L129").Value = 9684.599999999999
$ws.Range("M129").Value = 2503.21742
$ws.Range("N129").Value = -19684.6
$ws.Range("H132").Value = 1290.95

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 730.26666
$ws.Range("I97").Value = 664.0909
$ws.Range("J97").Value = 912.25
$ws.Range("K97").Value = 664.0909
$ws.Range("L97").Value = 912.25
$ws.Range("M97").Value = -168.0909
$ws.Range("N97").Value = -1904.25
$ws.Range("H132").Value = 3503.6938
$ws.Range("I132").Value = 3340.0625
$ws.Range("K132").Value = 10020.1875
$ws.Range("M132").Value = -7490.1875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3182.75
$ws.Range("I46").Value = 2238.5454
$ws.Range("J46").Value = 5260
$ws.Range("K46").Value = 2238.5454
$ws.Range("L46").Value = 5260
$ws.Range("M46").Value = -2050.5454
$ws.Range("N46").Value = -5636
$ws.Range("H55").Value = 419.0625
$ws.Range("I55").Value = 61.2
$ws.Range("J55").Value = 1015.5
$ws.Range("K55").Value = 61.2
$ws.Range("L55").Value = 1015.5
$ws.Range("M55").Value = 111.8
$ws.Range("N55").Value = -1361.5
$ws.Range("H132").Value = 3300.8552
$ws.Range("I132").Value = 2843.169
$ws.Range("K132").Value = 8529.507
$ws.Range("M132").Value = -5999.507
$ws.Range("H140").Value = 51557.5
$ws.Range("J140").Value = 51557.5
$ws.Range("L140").Value = 51557.5
$ws.Range("N140").Value = -61917.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 804602.9
$ws.Range("J81").Value = 4013.8333
$ws.Range("L81").Value = 8027.6666
$ws.Range("N81").Value = -10149.6666
$ws.Range("H84").Value = 804602.9
$ws.Range("J84").Value = 4013.8333
$ws.Range("L84").Value = 40138.333
$ws.Range("N84").Value = -50746.333
$ws.Range("H100").Value = 2156174
$ws.Range("I100").Value = 2463770
$ws.Range("K100").Value = 4927540
$ws.Range("M100").Value = -4926999
$ws.Range("H122").Value = 2919.8293
$ws.Range("I122").Value = 1490.7742
$ws.Range("J122").Value = 7349.9
$ws.Range("K122").Value = 4472.3226
$ws.Range("L122").Value = 22049.7
$ws.Range("M122").Value = -2022.3226
$ws.Range("N122").Value = -26949.7
$ws.Range("H132").Value = 2087.742
$ws.Range("I132").Value = 1694.4231
$ws.Range("J132").Value = 4133
$ws.Range("K132").Value = 5083.2693
$ws.Range("L132").Value = 12399
$ws.Range("M132").Value = -2553.2693
$ws.Range("N132").Value = -17459
$ws.Range("H138").Value = 125585.6
$ws.Range("J138").Value = 134557.25
$ws.Range("L138").Value = 134557.25
$ws.Range("N138").Value = -144837.25
